$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $vals) {
  foreach ($col in $vals.Keys) {
    $addr = "$col$row"
    $ws.Range($addr).Value = $vals[$col]
  }
}

# --- Rows 64/65, 66/67, 68/69: the match data in columns F..V got swapped
#     between each pair of rows (columns A-E, G untouched) ---

Set-Row 64 @{ F="Real Madrid"; H="Las Palmas"; I=0; J=1.2; K="23/09/2023 09:28"; L=1.18; M="27/09/2023 18:29"; N=6.76; O="23/09/2023 09:28"; P=8; Q="27/09/2023 18:29"; R=11.3; S="23/09/2023 09:28"; T=16.5; U="27/09/2023 18:29"; V="https://www.betexplorer.com/football/spain/laliga/real-madrid-las-palmas/GQHmRXXM/" }
Set-Row 65 @{ F="Ath Bilbao"; H="Getafe"; I=2; J=1.71; K="17/09/2023 09:02"; L=1.53; M="27/09/2023 18:31"; N=3.42; O="17/09/2023 09:02"; P=4.06; Q="27/09/2023 18:49"; R=5.44; S="17/09/2023 09:02"; T=7.73; U="27/09/2023 18:49"; V="https://www.betexplorer.com/football/spain/laliga/ath-bilbao-getafe/zgsFCYIT/" }

Set-Row 66 @{ F="Valencia"; H="Real Sociedad"; I=1; J=2.33; K="17/09/2023 09:02"; L=2.57; M="27/09/2023 21:27"; N=3.14; O="17/09/2023 09:02"; P=3.02; Q="27/09/2023 21:27"; R=3.26; S="17/09/2023 09:02"; T=3.31; U="27/09/2023 21:27"; V="https://www.betexplorer.com/football/spain/laliga/valencia-real-sociedad/M3IqSDIG/" }
Set-Row 67 @{ F="Cadiz CF"; H="Rayo Vallecano"; I=0; J=2.59; K="17/09/2023 09:02"; L=2.72; M="27/09/2023 21:19"; N=3.06; O="17/09/2023 09:02"; P=3.11; Q="27/09/2023 21:17"; R=3.11; S="17/09/2023 09:02"; T=3; U="27/09/2023 21:30"; V="https://www.betexplorer.com/football/spain/laliga/cadiz-rayo-vallecano/CEYt8hRp/" }

Set-Row 68 @{ F="Granada CF"; H="Betis"; I=1; J=3.25; K="23/09/2023 09:33"; L=2.79; M="28/09/2023 18:59"; N=3.38; O="23/09/2023 09:33"; P=3.46; Q="28/09/2023 18:57"; R=2.21; S="23/09/2023 09:33"; T=2.68; U="28/09/2023 18:57"; V="https://www.betexplorer.com/football/spain/laliga/granada-cf-betis/QyPzUZm4/" }
Set-Row 69 @{ F="Celta Vigo"; H="Alaves"; I=1; J=1.91; K="23/09/2023 09:33"; L=1.59; M="28/09/2023 18:58"; N=3.27; O="23/09/2023 09:33"; P=3.99; Q="28/09/2023 18:59"; R=4.39; S="23/09/2023 09:33"; T=6.84; U="28/09/2023 18:59"; V="https://www.betexplorer.com/football/spain/laliga/celta-vigo-alaves/fHGiQimT/" }

# --- New rows 127-129: three extra matches appended at the bottom ---

Set-Row 127 @{ A=126; B="spain"; C="laliga"; D="2023-2024"; E=45242.67708333334; F="Barcelona"; G=2; H="Alaves"; I=1; J=1.23; K="29/10/2023 11:02"; L=1.25; M="12/11/2023 16:07"; N=6.57; O="29/10/2023 11:02"; P=6.51; Q="12/11/2023 16:12"; R=12.29; S="29/10/2023 11:02"; T=12.09; U="12/11/2023 16:12"; V="https://www.betexplorer.com/football/spain/laliga/barcelona-alaves/8GjFNhSN/" }
Set-Row 128 @{ A=127; B="spain"; C="laliga"; D="2023-2024"; E=45242.77083333334; F="Sevilla"; G=1; H="Betis"; I=1; J=1.95; K="29/10/2023 11:02"; L=2.24; M="12/11/2023 18:22"; N=3.54; O="29/10/2023 11:02"; P=3.46; Q="12/11/2023 18:27"; R=4.07; S="29/10/2023 11:02"; T=3.42; U="12/11/2023 18:22"; V="https://www.betexplorer.com/football/spain/laliga/sevilla-betis/2D0fSGDh/" }
Set-Row 129 @{ A=128; B="spain"; C="laliga"; D="2023-2024"; E=45242.875; F="Atl. Madrid"; G=3; H="Villarreal"; I=1; J=1.58; K="29/10/2023 11:02"; L=1.34; M="12/11/2023 20:44"; N=4.3; O="29/10/2023 11:02"; P=5.81; Q="12/11/2023 20:59"; R=5.54; S="29/10/2023 11:02"; T=8.91; U="12/11/2023 20:59"; V="https://www.betexplorer.com/football/spain/laliga/atl-madrid-villarreal/fmcnUxqt/" }

# Match formatting used by the rest of the sheet: column A bold/centred with
# a thin border, column E using the existing custom date-time number format.
$idxRange = $ws.Range("A127:A129")
$idxRange.Font.Bold = $true
$idxRange.HorizontalAlignment = -4108
$idxRange.VerticalAlignment = -4160
$idxRange.Borders.LineStyle = 1

$dateRange = $ws.Range("E127:E129")
$dateRange.NumberFormat = "YYYY-MM-DD HH:MM:SS"
